# DDBJ BioSample "Virus; version 1.0" template: extend row 15 (the attribute/harmonized-name
# header row) with the remaining optional package attributes, columns M through AC.
#
# Each new header cell re-uses the formatting already applied to the existing optional ("yellow")
# header cells, e.g. C15 (isolation_source), and gets a cell comment with the attribute's
# definition, exactly like the existing A15:L15 headers already have.
#
# In particular, the harmonized name for the last of these columns (AC15) is written out directly
# as "temperature" (rather than the older, less descriptive "temp"), per the commit message:
# "temperature の harmonized name を temp から temperature に変更、分かりやすくするため".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cell to copy the "optional / yellow" header style from (e.g. C15 / isolation_source).
$templateRange = $ws.Range("C15")

# M15: altitude
$templateRange.Copy() | Out-Null
$ws.Range("M15").PasteSpecial(-4122) | Out-Null
$ws.Range("M15").Value = 'altitude'
$ws.Range("M15").AddComment('The altitude of the sample is the vertical distance between Earth''s surface above Sea Level and the sampled position in the air.') | Out-Null

# N15: biomaterial_provider
$templateRange.Copy() | Out-Null
$ws.Range("N15").PasteSpecial(-4122) | Out-Null
$ws.Range("N15").Value = 'biomaterial_provider'
$ws.Range("N15").AddComment('name and address of the lab or PI, or a culture collection identifier') | Out-Null

# O15: collected_by
$templateRange.Copy() | Out-Null
$ws.Range("O15").PasteSpecial(-4122) | Out-Null
$ws.Range("O15").Value = 'collected_by'
$ws.Range("O15").AddComment('Name of persons or institute who collected the sample') | Out-Null

# P15: culture_collection
$templateRange.Copy() | Out-Null
$ws.Range("P15").PasteSpecial(-4122) | Out-Null
$ws.Range("P15").Value = 'culture_collection'
$ws.Range("P15").AddComment('Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier') | Out-Null

# Q15: depth
$templateRange.Copy() | Out-Null
$ws.Range("Q15").PasteSpecial(-4122) | Out-Null
$ws.Range("Q15").Value = 'depth'
$ws.Range("Q15").AddComment('Depth is defined as the vertical distance below surface, e.g. for sediment or soil samples depth is measured from sediment or soil surface, respectively. Depth can be reported as an interval for subsurface samples.') | Out-Null

# R15: disease
$templateRange.Copy() | Out-Null
$ws.Range("R15").PasteSpecial(-4122) | Out-Null
$ws.Range("R15").Value = 'disease'
$ws.Range("R15").AddComment('list of diseases diagnosed; can include multiple diagnoses. the value of the field depends on host; for humans the terms should be chosen from DO (Disease Ontology), free text for non-human. For DO terms, please see http://gemina.svn.sourceforge.net/viewvc/gemina/trunk/Gemina/ontologies/gemina_symptom.obo?view=log') | Out-Null

# S15: env_biome
$templateRange.Copy() | Out-Null
$ws.Range("S15").PasteSpecial(-4122) | Out-Null
$ws.Range("S15").Value = 'env_biome'
$ws.Range("S15").AddComment('Descriptor of the broad ecological context of a sample. Examples include: desert, taiga or deciduous woodland. FAQ, http://trace.ddbj.nig.ac.jp/biosample/faq_e.html#biome-feature-material EnvO (v 2013-06-14) terms can be found via the link: http://www.environmentontology.org/Browse-EnvO') | Out-Null

# T15: genotype
$templateRange.Copy() | Out-Null
$ws.Range("T15").PasteSpecial(-4122) | Out-Null
$ws.Range("T15").Value = 'genotype'
$ws.Range("T15").AddComment('observed genotype') | Out-Null

# U15: host_tissue_sampled
$templateRange.Copy() | Out-Null
$ws.Range("U15").PasteSpecial(-4122) | Out-Null
$ws.Range("U15").Value = 'host_tissue_sampled'
$ws.Range("U15").AddComment('Type of tissue the initial sample was taken from. Controlled vocabulary, http://bioportal.bioontology.org/ontologies/1005') | Out-Null

# V15: identified_by
$templateRange.Copy() | Out-Null
$ws.Range("V15").PasteSpecial(-4122) | Out-Null
$ws.Range("V15").Value = 'identified_by'
$ws.Range("V15").AddComment('name of the taxonomist who identified the specimen') | Out-Null

# W15: lat_lon
$templateRange.Copy() | Out-Null
$ws.Range("W15").PasteSpecial(-4122) | Out-Null
$ws.Range("W15").Value = 'lat_lon'
$ws.Range("W15").AddComment('The geographical coordinates of the location where the sample was collected. Specify as decimal degrees latitude and longitude in format "d[d.dddd] N|S d[dd.dddd] W|E", eg, 47.94 N 28.12 W') | Out-Null

# X15: passage_history
$templateRange.Copy() | Out-Null
$ws.Range("X15").PasteSpecial(-4122) | Out-Null
$ws.Range("X15").Value = 'passage_history'
$ws.Range("X15").AddComment('Number of passages and passage method') | Out-Null

# Y15: samp_size
$templateRange.Copy() | Out-Null
$ws.Range("Y15").PasteSpecial(-4122) | Out-Null
$ws.Range("Y15").Value = 'samp_size'
$ws.Range("Y15").AddComment('Amount or size of sample (volume, mass or area) that was collected') | Out-Null

# Z15: serotype
$templateRange.Copy() | Out-Null
$ws.Range("Z15").PasteSpecial(-4122) | Out-Null
$ws.Range("Z15").Value = 'serotype'
$ws.Range("Z15").AddComment('Taxonomy below subspecies; a variety (in bacteria, fungi or virus) usually based on its antigenic properties. Same as serovar and serogroup. e.g. serotype="H1N1" in Influenza A virus CY098518.') | Out-Null

# AA15: specimen_voucher
$templateRange.Copy() | Out-Null
$ws.Range("AA15").PasteSpecial(-4122) | Out-Null
$ws.Range("AA15").Value = 'specimen_voucher'
$ws.Range("AA15").AddComment('Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a ''structured voucher''. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier') | Out-Null

# AB15: strain
$templateRange.Copy() | Out-Null
$ws.Range("AB15").PasteSpecial(-4122) | Out-Null
$ws.Range("AB15").Value = 'strain'
$ws.Range("AB15").AddComment('microbial or eukaryotic strain name') | Out-Null

# AC15: temperature
$templateRange.Copy() | Out-Null
$ws.Range("AC15").PasteSpecial(-4122) | Out-Null
$ws.Range("AC15").Value = 'temperature'
$ws.Range("AC15").AddComment('temperature of the sample at time of sampling') | Out-Null

